$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at the top, shifting existing data down by one row
$ws.Rows.Item(1).Insert()

# Set the new header cell content in B1 and make it bold
$ws.Range("B1").Value = "Elemente und Elementkombinationen"
$ws.Range("B1").Font.Bold = $true

# Update the active cell selection to B4
$ws.Range("B4").Select()
